# Remove workspace groups manual testing issue
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issues")

# Row 9 is "Manual Testing Workspace Groups" - delete the whole row,
# shifting everything below it up by one.
$ws.Rows.Item(9).Delete()
